# Update location tag BCR names: move the "BCR #" token from the end of
# the tag name to the beginning, e.g.
#   "Northern Pacific Rainforest (BCR 5)"  -> "BCR 5-Northern Pacific Rainforest"
#   "Northwestern Interior Forest(BCR 4)"  -> "BCR 4-Northwestern Interior Forest"
# Column C holds the slugified version of column B (spaces -> dashes, same
# casing as B), so it is regenerated from the new B value.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DefaultLocationTags")

for ($row = 98; $row -le 134; $row++) {
    $cell = $ws.Cells.Item($row, 2)
    $name = $cell.Value()

    # Pull the trailing "(BCR <number>)" (with or without a preceding space)
    # off of the name and rebuild as "BCR <number>-<name>".
    if ($name -match '^(.*?)\s*\(BCR\s*(\d+)\)\s*$') {
        $base = $matches[1]
        $num = $matches[2]
        $newName = "BCR $num-$base"

        $ws.Cells.Item($row, 2).Value = $newName
        $ws.Cells.Item($row, 3).Value = $newName.Replace(" ", "-")
    }
}

# Match the saved view state captured in the diff: scrolled so row 109 is at
# the top, with C139 as the active/selected cell.
$ws.Application.ActiveWindow.ScrollRow = 109
$ws.Range("C139").Select()
